$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3): micelle labels reused from shared strings (CTAB/TTAB/DTAB) ---
$ws.Range("H3").Value = "CTAB"
$ws.Range("I3").Value = "TTAB"
$ws.Range("J3").Value = "DTAB"

# --- New "solids" concentration column (G) ---
$ws.Range("G4").Value = 100
$ws.Range("G5").Value = 200
$ws.Range("G6").Value = 300

# --- New "ME" (Distância intermicelar) columns, filled as relative formulas ---
# each becomes an Excel "shared formula" across H4:H6 / I4:I6 / J4:J6
$ws.Range("H4:H6").Formula = "=D5"
$ws.Range("I4:I6").Formula = "=D10"
$ws.Range("J4:J6").Formula = "=D15"

# --- Number formats matching the target workbook ---
$ws.Range("H6").NumberFormat = "0.000"
$ws.Range("I4:I6").NumberFormat = "0.000"
$ws.Range("J4:J6").NumberFormat = "0.000"
$ws.Range("H4:H5").NumberFormat = "0.00"

# --- Final selection as left by the author ---
$ws.Range("J10").Select() | Out-Null
